# Evaluation workbook update: "Olivier's batch" on 25/05/2016.
# Adds a new shared string "25/05/2016 (Olivier)" and stamps it (with a
# green "filled-in" highlight) into the relevant Phase cells for several
# students, and also touches one pre-existing "18/05/2016 (Olivier)" cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = "25/05/2016 (Olivier)"
$oldText = "18/05/2016 (Olivier)"

# OLE (BGR) colors matching the workbook's "Good" palette:
#   font  FF006100 -> 24832
#   fill  FFC6EFCE -> 13561798
$fontColor = 24832
$fillColor = 13561798

function Set-Highlighted($addr) {
    $rng = $ws.Range($addr)
    $rng.Value = $newText
    $rng.Interior.Color = $fillColor
    $rng.Font.Color = $fontColor
}

function Set-Plain($addr, $text) {
    $rng = $ws.Range($addr)
    $rng.Value = $text
}

# Rows that receive the new highlighted value in one or more cells.
$highlighted = @(
    "B5",
    "C6", "D6", "E6", "F6",
    "B7",
    "B9",
    "B11",
    "B12",
    "B13", "C13",
    "B15",
    "B16",
    "C18", "D18", "E18",
    "B20",
    "B21", "C21", "D21",
    "B22",
    "B24",
    "B25",
    "B34",
    "B35",
    "B37",
    "B38",
    "B41",
    "B43",
    "B44",
    "B45", "C45", "D45",
    "B48",
    "B49",
    "B50",
    "B51",
    "B52",
    "B53"
)

foreach ($addr in $highlighted) {
    Set-Highlighted $addr
}

# Rows 28 and 39 get the new text but keep their plain (unfilled) style.
Set-Plain "B28" $newText
Set-Plain "B39" $newText
Set-Plain "C39" $newText

# Row 42 gets the highlight, but (matching the source data) references the
# existing "18/05/2016 (Olivier)" text rather than the new date.
Set-Highlighted "B42"
$ws.Range("B42").Value = $oldText

# View-state bookkeeping to match the saved window/scroll position.
$ws.Application.ActiveWindow.ScrollRow = 21
$ws.Range("B57").Select() | Out-Null
